# Refresh of the "2025" year-to-date poker standings (rows 201-208 of the
# zz_Poker_Yearly query table on Sheet1) following new data being pulled in
# from the underlying Power Query source (new month of results added).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 201 - Andy stays in 1st, points/chips/winnings/takehome increase
$ws.Range("D201").Value = 48
$ws.Range("F201").Value = 48
$ws.Range("G201").Value = 157700
$ws.Range("H201").Value = 220
$ws.Range("I201").Value = 130

# Row 202 - Prashant stays in 2nd, points/chips increase, takehome decreases
$ws.Range("D202").Value = 41
$ws.Range("F202").Value = 41
$ws.Range("G202").Value = 124000
$ws.Range("I202").Value = 50

# Row 203 - now Richard (was Maisy)
$ws.Range("B203").Value = "Richard"
$ws.Range("D203").Value = 33
$ws.Range("F203").Value = 33
$ws.Range("G203").Value = 102950
$ws.Range("H203").Value = 115
$ws.Range("I203").Value = 35
$ws.Range("K203").Value = 366

# Row 204 - now Pepe (was Richard)
$ws.Range("B204").Value = "Pepe"
$ws.Range("D204").Value = 32
$ws.Range("F204").Value = 32
$ws.Range("G204").Value = 95850
$ws.Range("H204").Value = 50
$ws.Range("I204").Value = -20
$ws.Range("K204").Value = 364

# Row 205 - now Maisy (was Matt)
$ws.Range("B205").Value = "Maisy"
$ws.Range("D205").Value = 31
$ws.Range("F205").Value = 31
$ws.Range("G205").Value = 89950
$ws.Range("H205").Value = 50
$ws.Range("K205").Value = 360

# Row 206 - now Mark (was Pepe)
$ws.Range("B206").Value = "Mark"
$ws.Range("D206").Value = 29
$ws.Range("F206").Value = 29
$ws.Range("G206").Value = 89150
$ws.Range("H206").Value = 65
$ws.Range("I206").Value = -25
$ws.Range("K206").Value = 361

# Row 207 - now Matt (was Mark)
$ws.Range("B207").Value = "Matt"
$ws.Range("D207").Value = 29
$ws.Range("F207").Value = 29
$ws.Range("G207").Value = 85200
$ws.Range("H207").Value = 60
$ws.Range("I207").Value = -20
$ws.Range("K207").Value = 362

# Row 208 - Anthony stays in 8th, points/chips/takehome change
$ws.Range("D208").Value = 26
$ws.Range("F208").Value = 26
$ws.Range("G208").Value = 66800
$ws.Range("I208").Value = -50
